# Regenerate Report for Handoff:
#  - New handoff markdown file id: f71bdb80-c406-4e89-9452-6702e929ad71
#                              ->  7cc8ff44-6ac0-4958-9a22-582003046462
#  - New xliff content hash:       8a10d5a4e56cff35bce71e86d3f308e0f8e7473d
#                              ->  5a34a9afc934c77e49db6a6e5385674c5835d76f
#  - Updated generation / handoff timestamps

$wb = $excel.ActiveWorkbook

$oldId = "f71bdb80-c406-4e89-9452-6702e929ad71"
$newId = "7cc8ff44-6ac0-4958-9a22-582003046462"

$oldHash = "8a10d5a4e56cff35bce71e86d3f308e0f8e7473d"
$newHash = "5a34a9afc934c77e49db6a6e5385674c5835d76f"

$newFileName   = "$newId.md"
$newPathName   = "e2e\$newId.md"
$newGenDate    = "2016-08-26 22:56:39"
$newZhXlf      = "$newId.$newHash.zh-cn.xlf"
$newZhDate     = "2016-08-26 22:56:35"
$newDeXlf      = "$newId.$newHash.de-de.xlf"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a62988451814d788ab3b6939447743bdcd1c0a2/e2e/$oldId.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newGenDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $githubBase, "", "", $newPathName)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $githubBase, "", "", $newFileName)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $githubBase, "", "", $newFileName)
